$d = $word.ActiveDocument

# --- Edit 1: "Tên đầy đủ" -> "Tên đầy đủ : " -------------------------------
# (adds a " :" run and a trailing " " run, both italic, after "đủ")
$rng = $d.Content
$found1 = $rng.Find.Execute("Tên đầy đủ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng.Collapse(0)
    $a0 = $rng.Start
    $rng.InsertAfter(" :")
    $rng.Collapse(0)
    $a1 = $rng.Start
    $rng.InsertAfter(" ")
    $rng.Collapse(0)
    $a2 = $rng.Start

    $d.Range($a0, $a1).Italic = 1
    $d.Range($a1, $a2).Italic = 1
}

# --- Edit 2: "Mô tả rủi ro:" -> "Mô tả rủi ro: tool upload chậm" -----------
# (adds a " tool upload " run and a "chậm" run, both italic)
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Mô tả rủi ro:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $b0 = $rng2.Start
    $rng2.InsertAfter(" tool upload ")
    $rng2.Collapse(0)
    $b1 = $rng2.Start
    $rng2.InsertAfter("chậm")
    $rng2.Collapse(0)
    $b2 = $rng2.Start

    $d.Range($b0, $b1).Italic = 1
    $d.Range($b1, $b2).Italic = 1
}
